# Automatic tracker update ("Actualización automática del tracker")
# Re-sync the match-tracking sheet: drop the stale duplicate Kalinskaya/Svitolina
# row (F=3.2), shift every following record up one row, refresh resultado/profit
# for the matches that have now finished, and append the newly-scheduled match
# (Popyrin vs Rune) at the bottom. This shrinks the used range from A1:H13 to
# A1:H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 5 (Anna Kalinskaya F=3.2, a stale duplicate of what is now row 8
# with F=3.4) is removed entirely; everything below it shifts up one row.
$ws.Range("A5:H5").Delete() | Out-Null

# Final data values for rows 2-11 (event_id, fecha, jugador_A, jugador_B,
# pronostico, cuota, resultado, profit) after the shift + result refresh.
$data = @(
    @(14265551, "2025-08-01", "Alejandro Davidovich Fokina", "Jakub Mensik", "Gana Jakub Mensik", 4.33, "Fallo", -1),
    @(14265545, "2025-08-01", "Andrey Rublev", "Lorenzo Sonego", "Gana Lorenzo Sonego", 3.2, "Fallo", -1),
    @(14265544, "2025-08-01", "Flavio Cobolli", "Fabian Marozsan", "Gana Flavio Cobolli", 2, "Acierto", 1),
    @(14266950, "2025-08-01", "Nishesh Basavareddy", "Zachary Svajda", "Gana Zachary Svajda", 5.5, "Acierto", 4.5),
    @(14266954, "2025-08-01", "Michael Zheng", "Yu Hsiou Hsu", "Gana Yu Hsiou Hsu", 2.38, "Acierto", 1.38),
    @(14267299, "2025-08-01", "Daniel Michalski", "Valentin Vacherot", "Gana Daniel Michalski", 3, "Acierto", 2),
    @(14266295, "2025-08-01", "Anna Kalinskaya", "Elina Svitolina", "Gana Anna Kalinskaya", 3.4, "Fallo", -1),
    @(14259084, "2025-08-02", "Abdullah Shelbayh", "Alexandr Binda", "Gana Alexandr Binda", 3.25, $null, $null),
    @(14266653, "2025-08-02", "Botic Van de Zandschulp", "Guy Den Ouden", "Gana Guy Den Ouden", 3.25, $null, $null),
    @(14266660, "2025-08-02", "Olle Wallin", "Vilius Gaubas", "Gana Olle Wallin", 4.5, $null, $null)
)

$ws.Range("B2:B12").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    if ($null -eq $row[6]) {
        $ws.Cells.Item($r, 7).Value = ""
    } else {
        $ws.Cells.Item($r, 7).Value = $row[6]
    }
    if ($null -eq $row[7]) {
        $ws.Cells.Item($r, 8).Value = ""
    } else {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
    $r++
}

# Append the newly-tracked match as the new last row (row 12).
$ws.Cells.Item(12, 1).Value = 14265591
$ws.Cells.Item(12, 2).Value = "2025-08-02"
$ws.Cells.Item(12, 3).Value = "Alexei Popyrin"
$ws.Cells.Item(12, 4).Value = "Holger Rune"
$ws.Cells.Item(12, 5).Value = "Gana Holger Rune"
$ws.Cells.Item(12, 6).Value = 1.4
$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(12, 8).Value = ""
